$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 80 (Leve Item ID 12605)
$ws.Range("H80").Value = 1096.5714
$ws.Range("J80").Value = 1332.6875
$ws.Range("L80").Value = 3998.0625
$ws.Range("N80").Value = -5994.0625
# Row 83 (Leve Item ID 12605)
$ws.Range("H83").Value = 1096.5714
$ws.Range("J83").Value = 1332.6875
$ws.Range("L83").Value = 11994.1875
$ws.Range("N83").Value = -21978.1875
# Row 99 (Leve Item ID 19883)
$ws.Range("H99").Value = 4585.25
$ws.Range("I99").Value = 1030.3334
$ws.Range("K99").Value = 3091.0002
$ws.Range("M99").Value = -1593.0002
# Row 127 (Leve Item ID 36114)
$ws.Range("H127").Value = 6281.2144
$ws.Range("I127").Value = 2483.111
$ws.Range("K127").Value = 7449.333
$ws.Range("M127").Value = -2489.333
# Row 134 (Leve Item ID 41997)
$ws.Range("H134").Value = 50000
$ws.Range("J134").Value = 50000
$ws.Range("L134").Value = 50000
$ws.Range("N134").Value = -60140
# Row 135 (Leve Item ID 44047)
$ws.Range("H135").Value = 5190.914
$ws.Range("I135").Value = 2030.2609
$ws.Range("J135").Value = 11248.833
$ws.Range("K135").Value = 18272.3481
$ws.Range("L135").Value = 101239.497
$ws.Range("M135").Value = -15737.3481
$ws.Range("N135").Value = -106309.497
# Row 136 (Leve Item ID 42164)
$ws.Range("H136").Value = 0
$ws.Range("J136").Value = 0
$ws.Range("L136").Value = 0
$ws.Range("N136").ClearContents()
# Row 137 (Leve Item ID 44013)
$ws.Range("H137").Value = 2400.5652
$ws.Range("I137").Value = 865.73334
$ws.Range("J137").Value = 5278.375
$ws.Range("K137").Value = 2597.20002
$ws.Range("L137").Value = 15835.125
$ws.Range("M137").Value = -47.20002000000022
$ws.Range("N137").Value = -20935.125
# Row 138 (Leve Item ID 44169)
$ws.Range("H138").Value = 2940.9204
$ws.Range("I138").Value = 1380.381
$ws.Range("K138").Value = 4141.143
$ws.Range("M138").Value = 998.857
# Row 140 (Leve Item ID 42459)
$ws.Range("H140").Value = 348748.75
$ws.Range("J140").Value = 341665
$ws.Range("L140").Value = 341665
$ws.Range("N140").Value = -352025

$ws = $wb.Worksheets.Item("ARM")
# Row 37 (Leve Item ID 3096)
$ws.Range("H37").Value = 6000
$ws.Range("J37").Value = 6000
$ws.Range("L37").Value = 6000
$ws.Range("N37").Value = -6546
# Row 61 (Leve Item ID 43999)
$ws.Range("H61").Value = 5063.6235
$ws.Range("I61").Value = 3200.6912
$ws.Range("K61").Value = 3200.6912
$ws.Range("M61").Value = -2988.6912
# Row 76 (Leve Item ID 10679)
$ws.Range("H76").Value = 46285.875
$ws.Range("J76").Value = 46285.875
$ws.Range("L76").Value = 46285.875
$ws.Range("N76").Value = -46961.875
# Row 79 (Leve Item ID 10679)
$ws.Range("H79").Value = 46285.875
$ws.Range("J79").Value = 46285.875
$ws.Range("L79").Value = 46285.875
$ws.Range("N79").Value = -48625.875
# Row 105 (Leve Item ID 18699)
$ws.Range("H105").Value = 20370
$ws.Range("J105").Value = 20370
$ws.Range("L105").Value = 20370
$ws.Range("N105").Value = -27358
# Row 132 (Leve Item ID 43997)
$ws.Range("H132").Value = 722472.0600000001
$ws.Range("I132").Value = 757550.1
$ws.Range("K132").Value = 2272650.3
$ws.Range("M132").Value = -2270120.3
# Row 136 (Leve Item ID 43999)
$ws.Range("H136").Value = 5063.6235
$ws.Range("I136").Value = 3200.6912
$ws.Range("K136").Value = 9602.0736
$ws.Range("M136").Value = -7052.0736

$ws = $wb.Worksheets.Item("BSM")
# Row 99 (Leve Item ID 19943)
$ws.Range("H99").Value = 9270.233
$ws.Range("I99").Value = 9058.378000000001
$ws.Range("K99").Value = 9058.378000000001
$ws.Range("M99").Value = -7560.378000000001

$ws = $wb.Worksheets.Item("CRP")
# Row 22 (Leve Item ID 5367)
$ws.Range("H22").Value = 894001.5
$ws.Range("I22").Value = 1190755.1
$ws.Range("K22").Value = 1190755.1
$ws.Range("M22").Value = -1190405.1
# Row 31 (Leve Item ID 44023)
$ws.Range("H31").Value = 8789.071
$ws.Range("I31").Value = 49792
$ws.Range("J31").Value = 4472.9736
$ws.Range("K31").Value = 49792
$ws.Range("L31").Value = 4472.9736
$ws.Range("M31").Value = -49497
$ws.Range("N31").Value = -5062.9736
# Row 33 (Leve Item ID 1836)
$ws.Range("H33").Value = 20017.5
$ws.Range("I33").Value = 10000
$ws.Range("J33").Value = 30035
$ws.Range("K33").Value = 10000
$ws.Range("L33").Value = 30035
$ws.Range("M33").Value = -9621
$ws.Range("N33").Value = -30793
# Row 34 (Leve Item ID 44023)
$ws.Range("H34").Value = 8789.071
$ws.Range("I34").Value = 49792
$ws.Range("J34").Value = 4472.9736
$ws.Range("K34").Value = 49792
$ws.Range("L34").Value = 4472.9736
$ws.Range("M34").Value = -49590
$ws.Range("N34").Value = -4876.9736
# Row 58 (Leve Item ID 44021)
$ws.Range("H58").Value = 9356.809999999999
$ws.Range("I58").Value = 7494
$ws.Range("J58").Value = 11050.272
$ws.Range("K58").Value = 7494
$ws.Range("L58").Value = 11050.272
$ws.Range("M58").Value = -7291
$ws.Range("N58").Value = -11456.272
# Row 122 (Leve Item ID 36196)
$ws.Range("H122").Value = 5323.2383
$ws.Range("I122").Value = 2333.0557
$ws.Range("J122").Value = 23264.334
$ws.Range("K122").Value = 6999.1671
$ws.Range("L122").Value = 69793.00199999999
$ws.Range("M122").Value = -4549.1671
$ws.Range("N122").Value = -74693.00199999999
# Row 136 (Leve Item ID 44021)
$ws.Range("H136").Value = 9356.809999999999
$ws.Range("I136").Value = 7494
$ws.Range("J136").Value = 11050.272
$ws.Range("K136").Value = 22482
$ws.Range("L136").Value = 33150.81600000001
$ws.Range("M136").Value = -19932
$ws.Range("N136").Value = -38250.81600000001

$ws = $wb.Worksheets.Item("CUL")
# Row 2 (Leve Item ID 4847)
$ws.Range("H2").Value = 80.47059
$ws.Range("I2").Value = 45.5
$ws.Range("J2").Value = 130.42857
$ws.Range("K2").Value = 273
$ws.Range("L2").Value = 782.57142
$ws.Range("M2").Value = -160
$ws.Range("N2").Value = -1008.57142
# Row 4 (Leve Item ID 4650)
$ws.Range("H4").Value = 5135851.5
$ws.Range("I4").Value = 5286904
$ws.Range("K4").Value = 15860712
$ws.Range("M4").Value = -15860600
# Row 5 (Leve Item ID 43974)
$ws.Range("H5").Value = 1025.85
$ws.Range("J5").Value = 1878.6154
$ws.Range("L5").Value = 5635.8462
$ws.Range("N5").Value = -5859.8462
# Row 18 (Leve Item ID 36056)
$ws.Range("H18").Value = 2479.3333
$ws.Range("I18").Value = 2560.6667
$ws.Range("K18").Value = 7682.000100000001
$ws.Range("M18").Value = -7513.000100000001
# Row 32 (Leve Item ID 4731)
$ws.Range("H32").Value = 837.6667
$ws.Range("I32").Value = 719.8570999999999
$ws.Range("J32").Value = 1250
$ws.Range("K32").Value = 2159.5713
$ws.Range("L32").Value = 3750
$ws.Range("M32").Value = -1876.5713
$ws.Range("N32").Value = -4316
# Row 46 (Leve Item ID 4701)
$ws.Range("H46").Value = 2633
$ws.Range("J46").Value = 2633
$ws.Range("L46").Value = 7899
$ws.Range("N46").Value = -8081
# Row 48 (Leve Item ID 4724)
$ws.Range("H48").Value = 0
$ws.Range("I48").Value = 0
$ws.Range("K48").Value = 0
$ws.Range("M48").ClearContents()
# Row 68 (Leve Item ID 12895)
$ws.Range("H68").Value = 1061.4286
# Row 71 (Leve Item ID 12895)
$ws.Range("H71").Value = 1061.4286
# Row 80 (Leve Item ID 12890)
$ws.Range("H80").Value = 2024.5
$ws.Range("J80").Value = 2998
$ws.Range("L80").Value = 8994
$ws.Range("N80").Value = -10866
# Row 83 (Leve Item ID 12890)
$ws.Range("H83").Value = 2024.5
$ws.Range("J83").Value = 2998
$ws.Range("L83").Value = 26982
$ws.Range("N83").Value = -36342
# Row 107 (Leve Item ID 27838)
$ws.Range("H107").Value = 2681.36
$ws.Range("J107").Value = 2881.8696
$ws.Range("L107").Value = 8645.6088
$ws.Range("N107").Value = -12485.6088
# Row 113 (Leve Item ID 27843)
$ws.Range("H113").Value = 2718422.8
$ws.Range("I113").Value = 10417254
$ws.Range("J113").Value = 1188.4117
$ws.Range("K113").Value = 31251762
$ws.Range("L113").Value = 3565.2351
$ws.Range("M113").Value = -31249592
$ws.Range("N113").Value = -7905.2351
# Row 119 (Leve Item ID 27873)
$ws.Range("H119").Value = 3974.25
$ws.Range("I119").Value = 3974.25
$ws.Range("K119").Value = 11922.75
$ws.Range("M119").Value = -7084.75
# Row 131 (Leve Item ID 36060)
$ws.Range("H131").Value = 6775.6665
$ws.Range("I131").Value = 7331.2
$ws.Range("J131").Value = 3998
$ws.Range("K131").Value = 21993.6
$ws.Range("L131").Value = 11994
$ws.Range("M131").Value = -16953.6
$ws.Range("N131").Value = -22074
# Row 135 (Leve Item ID 43974)
$ws.Range("H135").Value = 1025.85
$ws.Range("J135").Value = 1878.6154
$ws.Range("L135").Value = 16907.5386
$ws.Range("N135").Value = -21977.5386

$ws = $wb.Worksheets.Item("GSM")
# Row 92 (Leve Item ID 18094)
$ws.Range("H92").Value = 7995.5
$ws.Range("J92").Value = 7995.5
$ws.Range("L92").Value = 7995.5
$ws.Range("N92").Value = -11739.5
# Row 102 (Leve Item ID 36169)
$ws.Range("H102").Value = 4014.855
$ws.Range("I102").Value = 2778.6382
$ws.Range("J102").Value = 6655.864
$ws.Range("K102").Value = 2778.6382
$ws.Range("L102").Value = 6655.864
$ws.Range("M102").Value = -1156.6382
$ws.Range("N102").Value = -9899.864
# Row 122 (Leve Item ID 36182)
$ws.Range("H122").Value = 8560
$ws.Range("J122").Value = 19750
$ws.Range("L122").Value = 59250
$ws.Range("N122").Value = -64150
# Row 132 (Leve Item ID 44008)
$ws.Range("H132").Value = 4540.7754
$ws.Range("I132").Value = 4129.8823
$ws.Range("J132").Value = 5472.1333
$ws.Range("K132").Value = 12389.6469
$ws.Range("L132").Value = 16416.3999
$ws.Range("M132").Value = -9859.6469
$ws.Range("N132").Value = -21476.3999

$ws = $wb.Worksheets.Item("LTW")
# Row 101 (Leve Item ID 18549)
$ws.Range("H101").Value = 13095.714
$ws.Range("J101").Value = 13095.714
$ws.Range("L101").Value = 13095.714
$ws.Range("N101").Value = -19585.714
# Row 132 (Leve Item ID 44058)
$ws.Range("H132").Value = 5423.9644
$ws.Range("J132").Value = 8780
$ws.Range("L132").Value = 26340
$ws.Range("N132").Value = -31400
# Row 136 (Leve Item ID 44060)
$ws.Range("H136").Value = 7326.5835
$ws.Range("I136").Value = 7771.1333
$ws.Range("J136").Value = 6585.6665
$ws.Range("K136").Value = 23313.3999
$ws.Range("L136").Value = 19756.9995
$ws.Range("M136").Value = -20763.3999
$ws.Range("N136").Value = -24856.9995

$ws = $wb.Worksheets.Item("WVR")
# Row 126 (Leve Item ID 36210)
$ws.Range("H126").Value = 5613.1816
$ws.Range("I126").Value = 2714.4285
$ws.Range("J126").Value = 10686
$ws.Range("K126").Value = 8143.2855
$ws.Range("L126").Value = 32058
$ws.Range("M126").Value = -5673.2855
$ws.Range("N126").Value = -36998
# Row 132 (Leve Item ID 44029)
$ws.Range("H132").Value = 6234.9785
$ws.Range("I132").Value = 5660.9033
$ws.Range("J132").Value = 7347.25
$ws.Range("K132").Value = 16982.7099
$ws.Range("L132").Value = 22041.75
$ws.Range("M132").Value = -14452.7099
$ws.Range("N132").Value = -27101.75
# Row 136 (Leve Item ID 44031)
$ws.Range("H136").Value = 11915113
$ws.Range("I136").Value = 16676569
$ws.Range("J136").Value = 11472.417
$ws.Range("K136").Value = 50029707
$ws.Range("L136").Value = 34417.251
$ws.Range("M136").Value = -50027157
$ws.Range("N136").Value = -39517.251
